{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nfunction findIndex(text) {\n  for (let i = 0; i < items.length; i++) {\n    if (items[i].text === text) return i;\n  }\n  return -1;\n}\n\n// 1. Section heading: ERESOLVE fix -> React peer fix.\nconst headingIdx = findIndex(\"New in this update (Render npm ERESOLVE fix)\");\nif (headingIdx === -1) throw new Error(\"heading paragraph not found\");\nitems[headingIdx].insertText(\"New in this update (Render React peer fix)\", \"Replace\");\n\n// 2. Replace the four bullet/detail paragraphs describing the old\n//    date-fns/render.yaml fix with three paragraphs describing the new\n//    react/react-dom downgrade fix (net: one fewer paragraph).\nconst dateFnsIdx = findIndex(\"  - Changed `date-fns` from `^4.1.0` to `^3.6.0` in `frontend/package.json`.\");\nconst satisfiesIdx = findIndex(\"  - This satisfies `react-day-picker@8.10.1` peer requirement (`^2.28.0 || ^3.0.0`).\");\nconst renderYamlIdx = findIndex(\"- Updated `render.yaml` frontend build command back to:\");\nconst npmInstallIdx = findIndex(\"  - `npm install && npm run build`\");\nif (dateFnsIdx === -1 || satisfiesIdx === -1 || renderYamlIdx === -1 || npmInstallIdx === -1) {\n  throw new Error(\"detail paragraphs not found\");\n}\n\nitems[dateFnsIdx].insertText(\"  - `react-day-picker@8.10.1` supports React up to v18.\", \"Replace\");\nitems[satisfiesIdx].insertText(\n  \"  - Downgraded `react` and `react-dom` from `^19.0.0` to `^18.2.0` in `frontend/package.json`.\",\n  \"Replace\"\n);\nitems[renderYamlIdx].insertText(\n  \"- This resolves ERESOLVE error for `react-day-picker` peer dependency on Render.\",\n  \"Replace\"\n);\nitems[npmInstallIdx].delete();\n\n// 3. Git state section: bump last pushed commit hash and reword the\n//    \"fix is local\" note.\nconst commitIdx = findIndex(\"- Last pushed commit: 0369be1\");\nif (commitIdx === -1) throw new Error(\"commit paragraph not found\");\nitems[commitIdx].insertText(\"- Last pushed commit: 7d0fd23\", \"Replace\");\n\nconst localFixIdx = findIndex(\"- Current Render dependency fix is local and not pushed yet.\");\nif (localFixIdx === -1) throw new Error(\"local fix paragraph not found\");\nitems[localFixIdx].insertText(\"- Current React compatibility fix is local and not pushed yet.\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($p) {\n    $t = $p.Range.Text\n    if ($t.Length -gt 0 -and $t.EndsWith(\"`r\")) {\n        $t = $t.Substring(0, $t.Length - 1)\n    }\n    return $t\n}\n\nfunction Find-ParaByText($doc, $text) {\n    $count = $doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        if ((Get-ParaText $p) -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# 1. Section heading: ERESOLVE fix -> React peer fix.\n$heading = Find-ParaByText $d \"New in this update (Render npm ERESOLVE fix)\"\nif ($heading -eq $null) { throw \"heading paragraph not found\" }\n$heading.Range.Text = \"New in this update (Render React peer fix)\"\n\n# 2. Replace the four bullet/detail paragraphs describing the old\n#    date-fns/render.yaml fix with three paragraphs describing the new\n#    react/react-dom downgrade fix (net: one fewer paragraph). Locate all\n#    four by text before mutating anything, since edits/deletes shift\n#    paragraph indices.\n$dateFns = Find-ParaByText $d \"  - Changed ``date-fns`` from ``^4.1.0`` to ``^3.6.0`` in ``frontend/package.json``.\"\n$satisfies = Find-ParaByText $d \"  - This satisfies ``react-day-picker@8.10.1`` peer requirement (``^2.28.0 || ^3.0.0``).\"\n$renderYaml = Find-ParaByText $d \"- Updated ``render.yaml`` frontend build command back to:\"\n$npmInstall = Find-ParaByText $d \"  - ``npm install && npm run build``\"\nif ($dateFns -eq $null) { throw \"date-fns paragraph not found\" }\nif ($satisfies -eq $null) { throw \"satisfies paragraph not found\" }\nif ($renderYaml -eq $null) { throw \"render.yaml paragraph not found\" }\nif ($npmInstall -eq $null) { throw \"npm install paragraph not found\" }\n\n$dateFns.Range.Text = \"  - ``react-day-picker@8.10.1`` supports React up to v18.\"\n$satisfies.Range.Text = \"  - Downgraded ``react`` and ``react-dom`` from ``^19.0.0`` to ``^18.2.0`` in ``frontend/package.json``.\"\n$renderYaml.Range.Text = \"- This resolves ERESOLVE error for ``react-day-picker`` peer dependency on Render.\"\n$npmInstall.Range.Delete()\n\n# 3. Git state section: bump last pushed commit hash and reword the\n#    \"fix is local\" note.\n$commit = Find-ParaByText $d \"- Last pushed commit: 0369be1\"\nif ($commit -eq $null) { throw \"commit paragraph not found\" }\n$commit.Range.Text = \"- Last pushed commit: 7d0fd23\"\n\n$localFix = Find-ParaByText $d \"- Current Render dependency fix is local and not pushed yet.\"\nif ($localFix -eq $null) { throw \"local fix paragraph not found\" }\n$localFix.Range.Text = \"- Current React compatibility fix is local and not pushed yet.\"\n"}
